$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.879.29"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "2.319.91"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.89"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.30%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.446"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  -4.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.77%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.659.91"
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.104"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.839"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").Value = "2.321.13"
$ws.Range("E18").Value = "  +3.09%  "
$ws.Range("D19").Value = "43.785.89"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "0.0₃0981"
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.44%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "174.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.17%  "
$ws.Range("E31").Value = "  +4.61%  "
$ws.Range("E32").Value = "  -4.41%  "
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("E34").Value = "  +4.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0686"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.54%  "
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("E38").Value = "  -4.52%  "
$ws.Range("E39").Value = "  +3.46%  "
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0955"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("E49").Value = "  +4.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000209"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.78%  "
$ws.Range("D51").Value = "1.443.04"
$ws.Range("E51").Value = "  +0.12%  "
